$d = $word.ActiveDocument

# 1. Update the meeting date from "Febrero 16, 2015" to "Junio 17, 2015"
#    (this mirrors the Word UI edit that produced the recorded diff).
$null = $d.Content.Find.Execute(
    "Febrero 16, 2015", $true, $false, $false, $false, $false,
    $true, 1, $false, "Junio 17, 2015", 2)

# 2. Reproduce the run split / "_GoBack" bookmark placement left behind by
#    the editing session. Word marks the location of the most recent edit
#    with a hidden "_GoBack" bookmark (removing any previous one, since
#    bookmark names are unique), which also causes the surrounding run to
#    be split at that point.

# First carve off the trailing ", 2015" portion using a temporary bookmark
# so it becomes its own run (inserted after the "17").
$rTail = $d.Content
$null = $rTail.Find.Execute("Guadalajara, Jalisco  Junio 17")
$tailPos = $rTail.End
$tailRange = $d.Range($tailPos, $tailPos)
$d.Bookmarks.Add("TEMP_SPLIT", $tailRange)

# Now split off "Junio" from " 17" by placing "_GoBack" right after "Junio".
# Adding a bookmark named "_GoBack" also removes the pre-existing "_GoBack"
# bookmark (previously located near the end of the document) because
# bookmark names must be unique within the document.
$rMid = $d.Content
$null = $rMid.Find.Execute("Guadalajara, Jalisco  Junio")
$midPos = $rMid.End
$midRange = $d.Range($midPos, $midPos)
$d.Bookmarks.Add("_GoBack", $midRange)

# Remove the temporary helper bookmark; the run split it created persists.
$d.Bookmarks("TEMP_SPLIT").Delete()
